# CSSI 2019 budget updates (Amy Roberts)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Project period text (row 3) ---
$ws.Range("A3").Value = "Project Period:  11/1/2019-10/30/2022"

# --- Senior personnel: Gabriella Ramirez base salary (row 15) ---
$ws.Range("K15").Value = 43500

# --- Graduate Assistant(s) Enrolled staffing (row 21): 2 students @ 15 hrs/wk ---
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 15

# --- Undergraduate Assistant(s) Enrolled staffing (row 23): 15 hrs/wk ---
$ws.Range("L23").Value = 15

# --- Publication/Page Charges budget (row 42) raised to $1200/yr ---
$ws.Range("B42").Value = 1200
$ws.Range("C42").Value = 1200
$ws.Range("D42").Value = 1200

# --- Materials & Supplies (row 47) Year 1 reduced ---
$ws.Range("B47").Value = 2000

# --- New footnote re: computing devices as supplies (row 46) ---
$ws.Range("K46").Value = '" A computing device is considered a supply if the acquisition cost is less than the lesser of the capitalization level established by the proposer or $5,000, regardless of the length of its useful life. In the specific case of computing devices, charging as a direct cost is allowable for devices that are essential and allocable, but not solely dedicated, to the performance of the NSF project."'

# --- Permanent Equipment (row 58) Year 1 amount ---
$ws.Range("B58").Value = 12204

# --- Tuition row (row 65/66): new "cost per credit hour" label + credit hours ---
$ws.Range("L65").Value = "cost per credit hour"
$ws.Range("L66").Value = 373
$ws.Range("B66").Formula = "=2*K66*L66"

# --- Student Fees row (row 67) zeroed out ---
$ws.Range("B67").Value = 0

# --- Sheet view: scroll + selection to match reviewer's last position ---
$ws.Range("B62").Select() | Out-Null
